$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.950.89'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.539.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.74%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '611.05'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.47'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.538.32'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.78%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.488'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.54'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +9.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.435'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.82'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000218'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.139.52'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.542.17'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.978.03'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.18%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.61'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.76'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.95'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +11.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '451.23'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.637'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.90%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.14'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000130'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +5.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.684.27'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.89%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.04'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +11.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.29'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.69'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +12.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.53'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.171'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.23%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.92'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.74%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.25'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.90'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.534.63'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.11'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.32%  '

$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.33'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +8.94%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.10%  '

$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '177.30'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0909'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.57'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '31.25'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +14.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.889'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.22%  '

$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.97'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.79%  '

$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.33'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +9.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.57'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.71'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.259'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.65%  '
